# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (holding fund data, same layout as
#    2020-Q4 / 2021-Q2 / 2021-Q3) positioned right before the "总计" sheet.
# 2. Insert a new leading data row into "总计" for the 2022-Q1 totals and
#    renumber the existing rows' index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the "2022-Q1" worksheet with the fund-holding table
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q3")

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

# Move it so it sits immediately before "总计" (tab order: ... 2021-Q3,
# 2022-Q1, 总计).
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet.Move($totalSheet)
$ws = $wb.Worksheets.Item("2022-Q1")

# Copy the header-row formatting (bold + border + centered) and the
# first-column formatting from an existing quarter sheet so the new
# sheet's styles match the workbook's existing look exactly.
$refSheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$refSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("A2").Value = 0

# Fund code / numeric-looking figures must stay text (leading zeros,
# fixed decimal strings) - force text storage, then drop back to the
# default "Normal" style so no stray number-format style lingers on
# the cell (matches the plain, unstyled data cells elsewhere).
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "000049"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "中银标普全球精选自然资源等权重指数(QDII)"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.27"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "89.72"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1.33"
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.0036"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = 2

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: add the 2022-Q1 row at the top of the data
#    and renumber the existing rows.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()
$totals.Range("B2:D2").ClearFormats()

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
